$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder/relabel columns B-G ---
$ws.Range("B1").Value = "Parameters (B)"
$ws.Range("C1").Value = "Model Size (GB)"
$ws.Range("D1").Value = "Inference Time (s)"
$ws.Range("E1").Value = "Inference Memory (MiB)"
$ws.Range("F1").Value = "Winogrande (Accuracy %)"
$ws.Range("G1").Value = "Arc_challenge (Accuracy %)"

# --- Row 6 becomes "Flash Attention 2" or its associated metric row, with updated values ---
$ws.Range("A6").Value = "Flash Attention 2"
$ws.Range("B6").Value = 1.57
$ws.Range("C6").Value = 2.0699999999999998
$ws.Range("D6").Value = 10.975
$ws.Range("E6").Value = 1813.21
$ws.Range("F6").Value = 0.6369
$ws.Range("G6").Value = 0.4249

# --- Row 7 becomes a new "QLora" row with new data ---
$ws.Range("A7").Value = "QLora"
$ws.Range("B7").Value = 1.55
$ws.Range("C7").Value = 3.14
$ws.Range("D7").Value = 13.023999999999999
$ws.Range("E7").Value = 3340.33
$ws.Range("F7").Value = 0.49009000000000003
$ws.Range("G7").Value = 0.2201

# New row 7 (B7:G7) gets a distinct font color style (black RGB, explicit)
$ws.Range("B7:G7").Font.Color = 0

# --- Column widths for the two new metric columns ---
$ws.Columns.Item(6).ColumnWidth = 21.3
$ws.Columns.Item(7).ColumnWidth = 23.0

# --- View: zoom, window geometry + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 125
$win.Left = 39300
$win.Top = -7460
$win.Width = 21260
$win.Height = 10080
[void]$ws.Range("E6").Select()
